# Swap the full data (columns B, E..AD) between paired rows.
# Columns A (row index), C (Div) and D (Date) stay untouched because
# they are identical for the two rows in every pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(170,171),
    @(180,181),
    @(208,209),
    @(216,217),
    @(228,229)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Capture all values from both rows first, so the swap is not
    # affected by values already overwritten earlier in the loop.
    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range("$col$r1").Value()
        $vals2[$col] = $ws.Range("$col$r2").Value()
    }

    foreach ($col in $cols) {
        $ws.Range("$col$r1").Value = $vals2[$col]
        $ws.Range("$col$r2").Value = $vals1[$col]
    }
}
